$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. choices_female: add a "tuple" column (I) marking every row as a
#    single (non-coupled) entry, value 1 for every student.
# ---------------------------------------------------------------------
$wsFemale = $wb.Worksheets.Item("choices_female")

$wsFemale.Cells.Item(1, 9).Value = "tuple"
for ($r = 2; $r -le 57; $r++) {
    $wsFemale.Cells.Item($r, 9).Value = 1
}

$wsFemale.Range("I2:I57").Select()

# ---------------------------------------------------------------------
# 2. choices_male: couple the first two students ("A" and "B") into a
#    single row, tag every row with a "tuple" column that records how
#    many original students were folded into it (2 for the new
#    couple, 1 for everybody else).
# ---------------------------------------------------------------------
$wsMale = $wb.Worksheets.Item("choices_male")

# Remove the old "B" row - this shifts every later row up by one,
# exactly like the source data after the merge.
$wsMale.Rows(3).Delete()

# The surviving row 2 (originally "A") now represents the coupled pair.
$wsMale.Cells.Item(2, 1).Value = "A and B"

# Tag column with the tuple size.
$wsMale.Cells.Item(1, 7).Value = "tuple"
$wsMale.Cells.Item(2, 7).Value = 2
for ($r = 3; $r -le 35; $r++) {
    $wsMale.Cells.Item($r, 7).Value = 1
}

$wsMale.Range("G3").Select()

# choices_male becomes the active sheet/tab.
$wsMale.Activate()
